$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete existing data rows (2-8), keep header row 1
$ws.Range("A2:E8").Value = $null

# Make sure the date column keeps the original value as plain text (not parsed as a date)
$ws.Range("B2:B4").NumberFormat = "@"

# Row 2: Samsung Galaxy M30 complaint
$ws.Range("A2").Value = "Samsung Galaxy M30"
$ws.Range("B2").Value = "04/11/2019"
$ws.Range("C2").Value = "https://www.gsmarena.com/samsung_galaxy_m30-reviews-9505p1.php"
$ws.Range("D2").Value = "NFC,Fm Radio,"
$ws.Range("E2").Value = "M30 does not have gorilla glass, FM radio and NFC. A useless, worthless model.`nSamsung A9 2018 is far more better than M30 and A50`nA50 has gorilla glass 3 but does not have NFC and FM radio.`nSamsung is digging its own grave. Even Samsung A9 pro 2016 and other models have gorilla glass, FM radio and NFC"

# Row 3: Samsung Galaxy M20 - heating/update issue (reply)
$ws.Range("A3").Value = "Samsung Galaxy M20"
$ws.Range("B3").Value = "04/11/2019"
$ws.Range("C3").Value = "https://www.gsmarena.com/samsung_galaxy_m20-reviews-9506p1.php"
$ws.Range("D3").Value = "Video,Update,"
$ws.Range("E3").Value = "Anonymous, 9 hours agoHello im having  heating issu after the latest update and i chat with samsung support and they... moreWhy bro , they can't help? definitely they should help if your mobile phone doesn't have any physical damage. `nTry for second time (that time take a note about about the staff what they saying better you record a video and send email with the video to Samsung main head office - Korea regarding about the issue) they definitely help you, don't worry."

# Row 4: Samsung Galaxy M20 - original heating/update issue post
$ws.Range("A4").Value = "Samsung Galaxy M20"
$ws.Range("B4").Value = "04/11/2019"
$ws.Range("C4").Value = "https://www.gsmarena.com/samsung_galaxy_m20-reviews-9506p1.php"
$ws.Range("D4").Value = "Update,"
$ws.Range("E4").Value = "Hello im having  heating issu after the latest update and i chat with samsung support and they tell me go to the service center and the service center cant help me what can i do .`n"

# Restore the original (default) style/number format now that the text values are set
$ws.Range("B2:B4").Style = "Normal"
